$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1) Remove the obsolete "CQRules:CQBP-84--dependencies" rule row.
# ------------------------------------------------------------------
$depRow = $ws.Columns.Item(1).Find("CQRules:CQBP-84--dependencies")
if ($depRow -ne $null) {
    $ws.Rows.Item($depRow.Row).Delete() | Out-Null
}

# ------------------------------------------------------------------
# 2) "CloudServiceIncompatibleWorkflowProcess" changes severity from
#    Major to Blocker - remove the old (Major) row and re-insert it
#    directly after the last "Blocker" row (after "BannedPaths").
# ------------------------------------------------------------------
$cloudRow = $ws.Columns.Item(1).Find("CloudServiceIncompatibleWorkflowProcess")
if ($cloudRow -ne $null) {
    $ws.Rows.Item($cloudRow.Row).Delete() | Out-Null
}

$bannedRow = $ws.Columns.Item(1).Find("BannedPaths")
$insertAt = $bannedRow.Row + 1
$ws.Rows.Item($insertAt).Insert() | Out-Null
$ws.Cells.Item($insertAt, 1).Value = "CloudServiceIncompatibleWorkflowProcess"
$ws.Cells.Item($insertAt, 2).Value = "Usage of Cloud Service Incompatible Workflow Processes"
$ws.Cells.Item($insertAt, 3).Value = "Bug"
$ws.Cells.Item($insertAt, 4).Value = "Blocker"
$ws.Cells.Item($insertAt, 5).Value = "aem,cloud-service-compatibility"

# ------------------------------------------------------------------
# 3) Add the new "IndexDamAssetLucene" rule, right before
#    "ClientlibProxyResource" (end of the Bug/Minor block).
# ------------------------------------------------------------------
$clientlibRow = $ws.Columns.Item(1).Find("ClientlibProxyResource")
$newRow = $clientlibRow.Row
$ws.Rows.Item($newRow).Insert() | Out-Null
$ws.Cells.Item($newRow, 1).Value = "IndexDamAssetLucene"
$ws.Cells.Item($newRow, 2).Value = "Index customizations of the damAssetLucene Oak index should be properly structured."
$ws.Cells.Item($newRow, 3).Value = "Bug"
$ws.Cells.Item($newRow, 4).Value = "Minor"
$ws.Cells.Item($newRow, 5).Value = "aem,cloud-service-compatibility"

# ------------------------------------------------------------------
# 4) Mirror the final view state left by the author's editing session
#    (the last cell touched was the Tags cell of the new rule row).
# ------------------------------------------------------------------
$ws.Range("E" + $newRow).Select() | Out-Null
